$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Drop the stray "sa" run that precedes "<page>" at the very start
#    of the document (leftover scribble / typo, e.g. "sa<page>021r...").
# ---------------------------------------------------------------------
$d.Content.Find.Execute("sa<page>", $true, $false, $false, $false, $false, `
    $true, 1, $false, "<page>", 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Remove the two reviewer comments entirely (this also removes the
#    commentRangeStart/commentRangeEnd/commentReference markup from the
#    body - it leaves the plain "@" marker text behind, which we strip
#    in step 4 below).
# ---------------------------------------------------------------------
while ($d.Comments.Count -gt 0) {
    $d.Comments(1).Delete()
}

# ---------------------------------------------------------------------
# 3) "les assaillis @ reparent ..." -> "les assaillants ..."
#    First turn "i" (in "assaillis") into "ant" so it reads
#    "assaillants" (keeping that single letter's own run/formatting).
# ---------------------------------------------------------------------
$full = $d.Content.Text
$idx = $full.IndexOf("assaillis @")
$iStart = $idx + ("assaill").Length
$iEnd = $iStart + 1
$d.Range($iStart, $iEnd).Text = "ant"

# ---------------------------------------------------------------------
# 4) Delete the now-orphaned " @" (space + at-sign) that used to mark
#    the commented span, so the sentence reads
#    "les assaillants reparent Les assiegeants donnent une faulse".
# ---------------------------------------------------------------------
$full2 = $d.Content.Text
$spaceAtStart = $full2.IndexOf("assaillants @") + ("assaillants").Length
$spaceAtEnd = $spaceAtStart + (" @").Length
$d.Range($spaceAtStart, $spaceAtEnd).Delete()

Write-Output "edit complete"
